# Marca incidentes como "Resolvido" (coluna Status) e corrige o
# responsavel da linha 70 da aba ITI.

$wb = $excel.ActiveWorkbook

# --- Aba SPN ---
$wsSpn = $wb.Worksheets.Item("SPN")
$spnRows = @(76, 94, 98, 99, 102, 103)
foreach ($r in $spnRows) {
    $wsSpn.Range("I$r").Value = "Resolvido"
}

# --- Aba ITI ---
$wsIti = $wb.Worksheets.Item("ITI")

# Corrige o nome do responsavel na linha 70
$wsIti.Range("B70").Value = "Antônio Lucas"

$itiRows = @(12, 14, 15, 17, 22, 24, 28, 29, 30, 33, 34, 35, 37, 38, 40, 41, 42, 43, 44, 45, 46, 47, 48, 49, 50, 51, 53, 63, 70, 84, 88, 89, 91, 93, 98, 99, 100, 101, 102, 103, 104, 105, 106, 107, 108, 112, 113, 114, 115, 116, 117, 118, 119)
foreach ($r in $itiRows) {
    $wsIti.Range("I$r").Value = "Resolvido"
}
